$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$ws.Range("A3").Value  = "short 1000*1000"
$ws.Range("A4").Value  = "short 2000*2000"
$ws.Range("A5").Value  = "short 3000*3000"
$ws.Range("A6").Value  = "short 4000*4000"
$ws.Range("A7").Value  = "short 5000*5000"
$ws.Range("A8").Value  = "short 6000*6000"
$ws.Range("A9").Value  = "short 7000*7000"
$ws.Range("A10").Value = "short 8000*8000"
$ws.Range("A11").Value = "short 9000*9000"
$ws.Range("A12").Value = "short 10000*10000"

$ws.Range("R14").Select()
